$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows starting at row 16 (pushes existing rows 16-29 down to 20-33).
# The new rows inherit formatting (style) from the row above (row 15, style index 3).
$ws.Range("A16:A19").EntireRow.Insert()
$ws.Range("A16:A19").RowHeight = 17

# Row 15 used to read "Action: Generate a new customer when a user sets up billing for a student".
# That requirement is gone now - clear its text but keep the cell's style.
$ws.Range("A15").Value = ""

# New "requirements:" block describing the Subscription-signup feature.
$ws.Range("A16").Value = "requirements:"
$ws.Range("A17").Value = "If no Stripe Customer entity exists for user, render a button to sign user up."
$ws.Range("A18").Value = "Component with form to sign user up for a billing Subscription (create Subscription)"
$ws.Range("A19").Value = "Link existing ManageBilling Component to relevant dispatches"

# Row 21 (previously row 17, "4. Teachers can Create Subscription for student") is dropped entirely;
# clear its text but keep the style.
$ws.Range("A21").Value = ""

# Append the same three requirement lines again further down the sheet, reusing
# the same cell style/formatting as the "requirements:" block above.
$ws.Range("A17").Copy()
$ws.Range("A37:A39").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A37").Value = "If no Stripe Customer entity exists for user, render a button to sign user up."
$ws.Range("A38").Value = "Component with form to sign user up for a billing Subscription (create Subscription)"
$ws.Range("A39").Value = "Link existing ManageBilling Component to relevant dispatches"

# Update the active cell selection to match the authored file.
$ws.Range("F28").Select() | Out-Null
